$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# TOP ÁLBUNS: add a new "STREAMS" column (H) with header + 10 data values
# ---------------------------------------------------------------------------
$wsAlbuns = $wb.Worksheets.Item("TOP ÁLBUNS")

$wsAlbuns.Range("H1").Value = "STREAMS"
$wsAlbuns.Range("H1").Font.Bold = $true
$wsAlbuns.Range("H1").HorizontalAlignment = -4108
$wsAlbuns.Range("H1").NumberFormat = "0"

$streams = @(1421298139, 658516528, 2103224712, 656081062, 1979521945, 831781621, 5493041503, 5033709566, 5513097112, 413565042)
for ($i = 0; $i -lt $streams.Length; $i++) {
    $row = $i + 2
    $cell = $wsAlbuns.Cells.Item($row, 8)
    $cell.Value = $streams[$i]
    $cell.HorizontalAlignment = -4108
    $cell.NumberFormat = "0"
}

$wsAlbuns.Columns.Item(8).ColumnWidth = 13.0

# Update the selection on this sheet (it stays inactive / not the selected tab)
$wsAlbuns.Range("J8").Select()

# ---------------------------------------------------------------------------
# TOP ARTISTAS: keep selection, but it will no longer be the active tab
# ---------------------------------------------------------------------------
$wsArtistas = $wb.Worksheets.Item("TOP ARTISTAS")
$wsArtistas.Range("E3").Select()

# ---------------------------------------------------------------------------
# TOP MÚSICAS: becomes the active / selected tab, with a new selection
# ---------------------------------------------------------------------------
$wsMusicas = $wb.Worksheets.Item("TOP MÚSICAS")
$wsMusicas.Activate()
$wsMusicas.Range("C13").Select()
